$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.879.08"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "3.005.47"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "386.59"
$ws.Range("E5").Value = "  +3.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.14"
$ws.Range("E6").Value = "  +3.32%  "
$ws.Range("E7").Value = "  +1.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.599"
$ws.Range("E9").Value = "  +2.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.54"
$ws.Range("E10").Value = "  +2.03%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0850"
$ws.Range("E12").Value = "  +2.12%  "
$ws.Range("D13").Value = "3.482.56"
$ws.Range("E13").Value = "  +3.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.50"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("E15").Value = "  +4.23%  "
$ws.Range("D16").Value = "3.008.39"
$ws.Range("E16").Value = "  +3.74%  "
$ws.Range("E17").Value = "  +11.04%  "
$ws.Range("D18").Value = "51.855.48"
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.30"
$ws.Range("E19").Value = "  +2.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.48"
$ws.Range("E20").Value = "  +4.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.03"
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("D22").Value = "0.0₃0969"
$ws.Range("E22").Value = "  +3.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.35"
$ws.Range("E23").Value = "  +1.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.29"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("E25").Value = "  +9.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.41"
$ws.Range("E26").Value = "  +19.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.78"
$ws.Range("E27").Value = "  +24.68%  "
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("E29").Value = "  +13.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.17"
$ws.Range("E30").Value = "  +2.21%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.31"
$ws.Range("E33").Value = "  +3.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.13"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.08"
$ws.Range("E35").Value = "  -1.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0457"
$ws.Range("E36").Value = "  +8.78%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.06"
$ws.Range("E38").Value = "  +2.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.20"
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.62"
$ws.Range("E40").Value = "  +1.77%  "
$ws.Range("E41").Value = "  +1.42%  "
$ws.Range("E42").Value = "  +4.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "122.95"
$ws.Range("E43").Value = "  +3.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.97"
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.279"
$ws.Range("E45").Value = "  +18.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.06"
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.32"
$ws.Range("E47").Value = "  +5.96%  "
$ws.Range("E48").Value = "  +2.74%  "
$ws.Range("D49").Value = "2.047.26"
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0337"
$ws.Range("E50").Value = "  +9.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.872"
$ws.Range("E51").Value = "  +3.24%  "
